$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 157.8057217802531
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 415.5232050914462
